$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move the "Primary ROI" entries that were mistakenly placed in column C
# (rows 183-185, 187-192) back into column B. Row 186 (A186, "MB(L)")
# is a section header and stays where it is.
$rowsToMove = @(183, 184, 185, 187, 188, 189, 190, 191, 192)
foreach ($r in $rowsToMove) {
    $src = $ws.Cells.Item($r, 3)   # column C
    $dst = $ws.Cells.Item($r, 2)   # column B
    $src.Cut($dst)
}

# Append the new "NotPrimary" row at the bottom of the sheet.
$newRow = 231
$cell = $ws.Cells.Item($newRow, 1)
$cell.Value = "NotPrimary"
$cell.Style = $ws.Cells.Item(230, 1).Style

# Update the sheet's view/selection to match the new end of data.
$ws.Range("A231").Select()
$excel.ActiveWindow.ScrollRow = 218
